# Add a new "Discounted_Total" row below the existing "Total" row (A6:B6),
# computing 90% (i.e. a 10% discount) of the existing Total.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "Discounted_Total"
$ws.Range("B7").Formula = "=(1 - 0.1) * B6"

# Widen column A so the new, longer label fits (mirrors Excel's own
# auto-adjustment after typing a wider value into the column).
$ws.Columns("A:A").ColumnWidth = 15.417

# Leave the selection where Excel would land after entering data in B7.
$ws.Range("B8").Select() | Out-Null
